$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1013
$ws.Range("I32").Value = 576.6667
$ws.Range("J32").Value = 1106.5
$ws.Range("K32").Value = 576.6667
$ws.Range("L32").Value = 1106.5
$ws.Range("M32").Value = -250.6667
$ws.Range("N32").Value = -1758.5

$ws.Range("H106").Value = 3463.125
$ws.Range("I106").Value = 3252.5
$ws.Range("J106").Value = 3533.3333
$ws.Range("K106").Value = 3252.5
$ws.Range("L106").Value = 3533.3333
$ws.Range("M106").Value = -2621.5
$ws.Range("N106").Value = -4795.3333

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws.Range("H137").Value = 929.5128
$ws.Range("I137").Value = 762.1429000000001
$ws.Range("J137").Value = 1124.7778
$ws.Range("K137").Value = 2286.4287
$ws.Range("L137").Value = 3374.3334
$ws.Range("M137").Value = 263.5712999999996
$ws.Range("N137").Value = -8474.3334

$ws.Range("H140").Value = 84348
$ws.Range("J140").Value = 84348
$ws.Range("L140").Value = 84348
$ws.Range("N140").Value = -94708

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3661
$ws.Range("I61").Value = 3941.5715
$ws.Range("K61").Value = 3941.5715
$ws.Range("M61").Value = -3729.5715

$ws.Range("H74").Value = 2525.5518
$ws.Range("I74").Value = 2270.5881
$ws.Range("J74").Value = 2886.75
$ws.Range("K74").Value = 2270.5881
$ws.Range("L74").Value = 2886.75
$ws.Range("M74").Value = -1396.5881
$ws.Range("N74").Value = -4634.75

$ws.Range("H77").Value = 2525.5518
$ws.Range("I77").Value = 2270.5881
$ws.Range("J77").Value = 2886.75
$ws.Range("K77").Value = 11352.9405
$ws.Range("L77").Value = 14433.75
$ws.Range("M77").Value = -6984.940500000001
$ws.Range("N77").Value = -23169.75

$ws.Range("H110").Value = 1199.3448
$ws.Range("I110").Value = 1009.63635
$ws.Range("J110").Value = 1795.5714
$ws.Range("K110").Value = 1009.63635
$ws.Range("L110").Value = 1795.5714
$ws.Range("M110").Value = 1035.36365
$ws.Range("N110").Value = -5885.5714

$ws.Range("H133").Value = 75065.25
$ws.Range("J133").Value = 75065.25
$ws.Range("L133").Value = 75065.25
$ws.Range("N133").Value = -80125.25

$ws.Range("H136").Value = 3661
$ws.Range("I136").Value = 3941.5715
$ws.Range("K136").Value = 11824.7145
$ws.Range("M136").Value = -9274.7145

$ws.Range("H138").Value = 61196.668
$ws.Range("J138").Value = 61196.668
$ws.Range("L138").Value = 61196.668
$ws.Range("N138").Value = -71476.66800000001

$ws.Range("H139").Value = 64350
$ws.Range("J139").Value = 64350
$ws.Range("L139").Value = 64350
$ws.Range("N139").Value = -74630

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2396
$ws.Range("I99").Value = 1693.3334
$ws.Range("K99").Value = 1693.3334
$ws.Range("M99").Value = -195.3334

$ws.Range("H105").Value = 1693.4138
$ws.Range("I105").Value = 1692.36
$ws.Range("K105").Value = 1692.36
$ws.Range("M105").Value = 54.6400000000001

$ws.Range("H132").Value = 50780
$ws.Range("J132").Value = 50780
$ws.Range("L132").Value = 50780
$ws.Range("N132").Value = -60900

$ws.Range("H134").Value = 3912.276
$ws.Range("I134").Value = 1090.8975
$ws.Range("J134").Value = 9703.526
$ws.Range("K134").Value = 3272.6925
$ws.Range("L134").Value = 29110.578
$ws.Range("M134").Value = -737.6925000000001
$ws.Range("N134").Value = -34180.578

$ws.Range("H138").Value = 45776.668
$ws.Range("J138").Value = 45776.668
$ws.Range("L138").Value = 45776.668
$ws.Range("N138").Value = -56056.668

$ws.Range("H140").Value = 70860
$ws.Range("J140").Value = 70860
$ws.Range("L140").Value = 70860
$ws.Range("N140").Value = -81220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3881.3835
$ws.Range("I31").Value = 3206.9395
$ws.Range("J31").Value = 4437.8
$ws.Range("K31").Value = 3206.9395
$ws.Range("L31").Value = 4437.8
$ws.Range("M31").Value = -2911.9395
$ws.Range("N31").Value = -5027.8

$ws.Range("H34").Value = 3881.3835
$ws.Range("I34").Value = 3206.9395
$ws.Range("J34").Value = 4437.8
$ws.Range("K34").Value = 3206.9395
$ws.Range("L34").Value = 4437.8
$ws.Range("M34").Value = -3004.9395
$ws.Range("N34").Value = -4841.8

$ws.Range("H58").Value = 1215.3334
$ws.Range("I58").Value = 721.7143
$ws.Range("J58").Value = 1462.1428
$ws.Range("K58").Value = 721.7143
$ws.Range("L58").Value = 1462.1428
$ws.Range("M58").Value = -518.7143
$ws.Range("N58").Value = -1868.1428

$ws.Range("H105").Value = 1011.6667
$ws.Range("I105").Value = 664
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 664
$ws.Range("L105").Value = 2750
$ws.Range("M105").Value = 1083
$ws.Range("N105").Value = -6244

$ws.Range("H136").Value = 1215.3334
$ws.Range("I136").Value = 721.7143
$ws.Range("J136").Value = 1462.1428
$ws.Range("K136").Value = 2165.1429
$ws.Range("L136").Value = 4386.428400000001
$ws.Range("M136").Value = 384.8571000000002
$ws.Range("N136").Value = -9486.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 133.10527
$ws.Range("I12").Value = 37
$ws.Range("K12").Value = 111
$ws.Range("M12").Value = 62

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 35645
$ws.Range("J133").Value = 35645
$ws.Range("L133").Value = 35645
$ws.Range("N133").Value = -45765

$ws.Range("H138").Value = 64110
$ws.Range("J138").Value = 64110
$ws.Range("L138").Value = 64110
$ws.Range("N138").Value = -74390

$ws.Range("H140").Value = 99734.5
$ws.Range("J140").Value = 99734.5
$ws.Range("L140").Value = 99734.5
$ws.Range("N140").Value = -110094.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 958.8570999999999
$ws.Range("I55").Value = 1477.7778
$ws.Range("J55").Value = 569.6667
$ws.Range("K55").Value = 1477.7778
$ws.Range("L55").Value = 569.6667
$ws.Range("M55").Value = -1304.7778
$ws.Range("N55").Value = -915.6667

$ws.Range("H122").Value = 2711.7693
$ws.Range("I122").Value = 2449.2104
$ws.Range("J122").Value = 3424.4285
$ws.Range("K122").Value = 7347.6312
$ws.Range("L122").Value = 10273.2855
$ws.Range("M122").Value = -4897.6312
$ws.Range("N122").Value = -15173.2855

$ws.Range("H133").Value = 80483.88
$ws.Range("J133").Value = 80483.88
$ws.Range("L133").Value = 80483.88
$ws.Range("N133").Value = -85543.88

$ws.Range("H139").Value = 55920
$ws.Range("J139").Value = 55920
$ws.Range("L139").Value = 55920
$ws.Range("N139").Value = -66200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 20000000
$ws.Range("I29").Value = 20000000
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 20000000
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -19999710

$ws.Range("H136").Value = 978.7179599999999
$ws.Range("I136").Value = 825.0476
$ws.Range("J136").Value = 1158
$ws.Range("K136").Value = 2475.1428
$ws.Range("L136").Value = 3474
$ws.Range("M136").Value = 74.85719999999992
$ws.Range("N136").Value = -8574

$ws.Range("H138").Value = 47344.445
$ws.Range("J138").Value = 47344.445
$ws.Range("L138").Value = 47344.445
$ws.Range("N138").Value = -57624.445

$ws.Range("H140").Value = 19464.5
$ws.Range("J140").Value = 19464.5
$ws.Range("L140").Value = 19464.5
$ws.Range("N140").Value = -29824.5

$ws.Range("H141").Value = 57612.5
$ws.Range("J141").Value = 58985.715
$ws.Range("L141").Value = 58985.715
$ws.Range("N141").Value = -69345.715
